$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "86+2="
$t.Cell(1,2).Range.Text = "87-21="
$t.Cell(1,3).Range.Text = "91-71="
$t.Cell(1,4).Range.Text = "96-50="
$t.Cell(1,5).Range.Text = "68-53="
$t.Cell(2,1).Range.Text = "5+39="
$t.Cell(2,2).Range.Text = "28+25="
$t.Cell(2,3).Range.Text = "43+1="
$t.Cell(2,4).Range.Text = "15+52="
$t.Cell(2,5).Range.Text = "37-3="
$t.Cell(3,1).Range.Text = "29-15="
$t.Cell(3,2).Range.Text = "64+9="
$t.Cell(3,3).Range.Text = "69-13="
$t.Cell(3,4).Range.Text = "26-0="
$t.Cell(3,5).Range.Text = "88-6="
$t.Cell(4,1).Range.Text = "83-53="
$t.Cell(4,2).Range.Text = "67+5="
$t.Cell(4,3).Range.Text = "25+6="
$t.Cell(4,4).Range.Text = "67+20="
$t.Cell(4,5).Range.Text = "82-1="
$t.Cell(5,1).Range.Text = "36+46="
$t.Cell(5,2).Range.Text = "84-65="
$t.Cell(5,3).Range.Text = "2+96="
$t.Cell(5,4).Range.Text = "31-22="
$t.Cell(5,5).Range.Text = "0+58="
$t.Cell(6,1).Range.Text = "3+2="
$t.Cell(6,2).Range.Text = "70-10="
$t.Cell(6,3).Range.Text = "68-65="
$t.Cell(6,4).Range.Text = "21+41="
$t.Cell(6,5).Range.Text = "77-49="
$t.Cell(7,1).Range.Text = "23-6="
$t.Cell(7,2).Range.Text = "73-55="
$t.Cell(7,3).Range.Text = "38-37="
$t.Cell(7,4).Range.Text = "23-5="
$t.Cell(7,5).Range.Text = "60-49="
$t.Cell(8,1).Range.Text = "61-8="
$t.Cell(8,2).Range.Text = "46+36="
$t.Cell(8,3).Range.Text = "71-32="
$t.Cell(8,4).Range.Text = "88-36="
$t.Cell(8,5).Range.Text = "57+10="
$t.Cell(9,1).Range.Text = "62-35="
$t.Cell(9,2).Range.Text = "16+12="
$t.Cell(9,3).Range.Text = "46+44="
$t.Cell(9,4).Range.Text = "91-31="
$t.Cell(9,5).Range.Text = "19+60="
$t.Cell(10,1).Range.Text = "23+29="
$t.Cell(10,2).Range.Text = "16+65="
$t.Cell(10,3).Range.Text = "89-37="
$t.Cell(10,4).Range.Text = "61-42="
$t.Cell(10,5).Range.Text = "7-5="
$t.Cell(11,1).Range.Text = "20+79="
$t.Cell(11,2).Range.Text = "76-13="
$t.Cell(11,3).Range.Text = "10+79="
$t.Cell(11,4).Range.Text = "18+14="
$t.Cell(11,5).Range.Text = "36-8="
$t.Cell(12,1).Range.Text = "77-73="
$t.Cell(12,2).Range.Text = "78-7="
$t.Cell(12,3).Range.Text = "53+27="
$t.Cell(12,4).Range.Text = "66-35="
$t.Cell(12,5).Range.Text = "23-8="
$t.Cell(13,1).Range.Text = "48-36="
$t.Cell(13,2).Range.Text = "82+3="
$t.Cell(13,3).Range.Text = "44+31="
$t.Cell(13,4).Range.Text = "70+13="
$t.Cell(13,5).Range.Text = "25-6="
$t.Cell(14,1).Range.Text = "74-6="
$t.Cell(14,2).Range.Text = "2+47="
$t.Cell(14,3).Range.Text = "98-95="
$t.Cell(14,4).Range.Text = "0+43="
$t.Cell(14,5).Range.Text = "80-60="
$t.Cell(15,1).Range.Text = "56-3="
$t.Cell(15,2).Range.Text = "79-56="
$t.Cell(15,3).Range.Text = "23-23="
$t.Cell(15,4).Range.Text = "35-26="
$t.Cell(15,5).Range.Text = "58-15="
$t.Cell(16,1).Range.Text = "42-24="
$t.Cell(16,2).Range.Text = "5+13="
$t.Cell(16,3).Range.Text = "64-54="
$t.Cell(16,4).Range.Text = "9-4="
$t.Cell(16,5).Range.Text = "92-18="
$t.Cell(17,1).Range.Text = "29+65="
$t.Cell(17,2).Range.Text = "29-0="
$t.Cell(17,3).Range.Text = "8+46="
$t.Cell(17,4).Range.Text = "32+1="
$t.Cell(17,5).Range.Text = "46+30="
$t.Cell(18,1).Range.Text = "20+30="
$t.Cell(18,2).Range.Text = "39+3="
$t.Cell(18,3).Range.Text = "95-82="
$t.Cell(18,4).Range.Text = "16+2="
$t.Cell(18,5).Range.Text = "26+1="
$t.Cell(19,1).Range.Text = "78+20="
$t.Cell(19,2).Range.Text = "91-36="
$t.Cell(19,3).Range.Text = "95-83="
$t.Cell(19,4).Range.Text = "56+11="
$t.Cell(19,5).Range.Text = "22+32="
$t.Cell(20,1).Range.Text = "20-9="
$t.Cell(20,2).Range.Text = "49-38="
$t.Cell(20,3).Range.Text = "33+43="
$t.Cell(20,4).Range.Text = "31+9="
$t.Cell(20,5).Range.Text = "16+7="
